$d = $word.ActiveDocument

function Find-ParagraphByPrefix($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# --- 1. Insert "*New train asset" paragraph after "Add train door + Train SFX" ---
$pTrain = Find-ParagraphByPrefix("Add train door")
$pTrain.Range.InsertParagraphAfter()
$pNewTrain = $pTrain.Next()
$pNewTrain.Range.Text = "*New train asset"

# --- 2. Merge "Add signal to teach the player how to jump" + " - first jump room" runs ---
$dash = [char]8211
$pSignal = Find-ParagraphByPrefix("Add signal to teach")
$rSignal = $pSignal.Range
$mergedText = "Add signal to teach the player how to jump" + " " + $dash + " first jump room"
$rSignal.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2) | Out-Null

# --- 3. Insert a new run "*" before "- UI's" (kept as a distinct run) ---
$pUi = Find-ParagraphByPrefix("- UI")
$rUi = $pUi.Range
$insStar = $d.Range($rUi.Start, $rUi.Start)
$insStar.InsertBefore("*")
$rStar = $d.Range($rUi.Start, $rUi.Start + 1)
$rStar.Font.Bold = $true
$rStar.Font.Bold = $false

# --- 4. Strike-through the two bug paragraphs ---
$pStream = Find-ParagraphByPrefix("- Fix Level Streaming")
$pStream.Range.Font.StrikeThrough = $true
$pCheckpoint = Find-ParagraphByPrefix("- Fix Checkpoint")
$pCheckpoint.Range.Font.StrikeThrough = $true

Write-Output "Done"
